$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "descr"
$ws.Range("D1").Value = "is_active"

# --- Data rows (rows 2-9) ---
$data = @(
    @("eng", "txt",  "Text File",     $true),
    @("eng", "xml",  "XML File",      $true),
    @("eng", "json", "Json File",     $true),
    @("fra", "txt",  "Fichier texte", $true),
    @("fra", "xml",  "Fichier XML",   $true),
    @("fra", "json", "Fichier Json",  $true),
    @("eng", "html", "html file",     $true),
    @("fra", "html", "Fichier html",  $true)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]

    $row++
}

# Column A (rows 2-9) gets the same bordered/bold/centered style as the
# header row - copy the format from A1 down via paste-special (formats only).
$ws.Range("A1").Copy()
$ws.Range("A2:A9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
